$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 18: You Grow, Girl
$ws.Range("H18").Value = 328
$ws.Range("I18").Value = 285
$ws.Range("J18").Value = 500
$ws.Range("K18").Value = 285
$ws.Range("L18").Value = 500
$ws.Range("M18").Value = -1
$ws.Range("N18").Value = -1068

# Row 31: Hush Little Wailer
$ws.Range("H31").Value = 1325
$ws.Range("I31").Value = 500
$ws.Range("J31").Value = 2150
$ws.Range("K31").Value = 1500
$ws.Range("L31").Value = 6450
$ws.Range("M31").Value = -1270
$ws.Range("N31").Value = -6910

# Row 43: Growing Is Knowing
$ws.Range("H43").Value = 1038.3
$ws.Range("I43").Value = 1000.3333
$ws.Range("J43").Value = 1054.5714
$ws.Range("K43").Value = 1000.3333
$ws.Range("L43").Value = 1054.5714
$ws.Range("M43").Value = -931.3333
$ws.Range("N43").Value = -1192.5714

# Row 51: A Bile Business
$ws.Range("H51").Value = 1359.8
$ws.Range("I51").Value = 1049.5
$ws.Range("J51").Value = 1566.6666
$ws.Range("K51").Value = 1049.5
$ws.Range("L51").Value = 1566.6666
$ws.Range("M51").Value = -565.5
$ws.Range("N51").Value = -2534.6666

# Row 86: Filling in the Blanks
$ws.Range("H86").Value = 66697996
$ws.Range("I86").Value = 84249380
$ws.Range("J86").Value = 2760.8
$ws.Range("K86").Value = 84249380
$ws.Range("L86").Value = 2760.8
$ws.Range("M86").Value = -84248257
$ws.Range("N86").Value = -5006.8

# Row 88: The Grave of Hemlock Groves
$ws.Range("H88").Value = 31500.4
$ws.Range("I88").Value = 500
$ws.Range("J88").Value = 39250.5
$ws.Range("K88").Value = 500
$ws.Range("L88").Value = 39250.5
$ws.Range("M88").Value = -94
$ws.Range("N88").Value = -40062.5

# Row 89: Ink into Antiquity (L)
$ws.Range("H89").Value = 66697996
$ws.Range("I89").Value = 84249380
$ws.Range("J89").Value = 2760.8
$ws.Range("K89").Value = 421246900
$ws.Range("L89").Value = 13804
$ws.Range("M89").Value = -421241284
$ws.Range("N89").Value = -25036

# Row 91: Dappling the Highlands (L)
$ws.Range("H91").Value = 31500.4
$ws.Range("I91").Value = 500
$ws.Range("J91").Value = 39250.5
$ws.Range("K91").Value = 500
$ws.Range("L91").Value = 39250.5
$ws.Range("M91").Value = 904
$ws.Range("N91").Value = -42058.5

# Row 103: Let Loose the Juice
$ws.Range("H103").Value = 125372.75
$ws.Range("I103").Value = 143243.14
$ws.Range("J103").Value = 280
$ws.Range("K103").Value = 429729.42
$ws.Range("L103").Value = 840
$ws.Range("M103").Value = -429143.42
$ws.Range("N103").Value = -2012

# Row 116: Growing Up
$ws.Range("H116").Value = 5835.577
$ws.Range("I116").Value = 7143.1055
$ws.Range("J116").Value = 2286.5715
$ws.Range("K116").Value = 7143.1055
$ws.Range("L116").Value = 2286.5715
$ws.Range("M116").Value = -3701.1055
$ws.Range("N116").Value = -9170.5715

# Row 132: Fast-forwarding Flora
$ws.Range("H132").Value = 3834.6
$ws.Range("I132").Value = 3622.0344
$ws.Range("J132").Value = 9999
$ws.Range("K132").Value = 10866.1032
$ws.Range("L132").Value = 29997
$ws.Range("M132").Value = -8336.1032
$ws.Range("N132").Value = -35057

$ws = $wb.Worksheets.Item("ARM")
# Row 103: Sweeping the Legs
$ws.Range("H103").Value = 100000
$ws.Range("J103").Value = 100000
$ws.Range("L103").Value = 100000
$ws.Range("N103").Value = -102344

# Row 132: Don't Bore Me, Ore Me
$ws.Range("H132").Value = 3632.1904
$ws.Range("I132").Value = 3312.6296
$ws.Range("J132").Value = 4207.4
$ws.Range("K132").Value = 9937.888800000001
$ws.Range("L132").Value = 12622.2
$ws.Range("M132").Value = -7407.888800000001
$ws.Range("N132").Value = -17682.2

$ws = $wb.Worksheets.Item("BSM")
# Row 99: Meddle in Metal
$ws.Range("H99").Value = 916.6667
$ws.Range("I99").Value = 800
$ws.Range("J99").Value = 1220
$ws.Range("K99").Value = 800
$ws.Range("L99").Value = 1220
$ws.Range("M99").Value = 698
$ws.Range("N99").Value = -4216

# Row 134: Ruthenium Supremium
$ws.Range("H134").Value = 2217.1482
$ws.Range("I134").Value = 2032.0476
$ws.Range("J134").Value = 2865
$ws.Range("K134").Value = 6096.142800000001
$ws.Range("L134").Value = 8595
$ws.Range("M134").Value = -3561.142800000001
$ws.Range("N134").Value = -13665

$ws = $wb.Worksheets.Item("CRP")
# Row 31: Wall Not Found
$ws.Range("H31").Value = 2314
$ws.Range("I31").Value = 1898.9166
$ws.Range("J31").Value = 2867.4443
$ws.Range("K31").Value = 1898.9166
$ws.Range("L31").Value = 2867.4443
$ws.Range("M31").Value = -1603.9166
$ws.Range("N31").Value = -3457.4443

# Row 34: Armoires of the Rich and Famous
$ws.Range("H34").Value = 2314
$ws.Range("I34").Value = 1898.9166
$ws.Range("J34").Value = 2867.4443
$ws.Range("K34").Value = 1898.9166
$ws.Range("L34").Value = 2867.4443
$ws.Range("M34").Value = -1696.9166
$ws.Range("N34").Value = -3271.4443

# Row 132: Hull Lotta Damage
$ws.Range("H132").Value = 9261652
$ws.Range("I132").Value = 2897.8
$ws.Range("J132").Value = 12822712
$ws.Range("K132").Value = 8693.400000000001
$ws.Range("L132").Value = 38468136
$ws.Range("M132").Value = -6163.400000000001
$ws.Range("N132").Value = -38473196

$ws = $wb.Worksheets.Item("CUL")
# Row 2: Pork Is a Salty Food
$ws.Range("H2").Value = 275.75
$ws.Range("I2").Value = 201
$ws.Range("J2").Value = 300.66666
$ws.Range("K2").Value = 1206
$ws.Range("L2").Value = 1803.99996
$ws.Range("M2").Value = -1093
$ws.Range("N2").Value = -2029.99996

# Row 12: Butter Me Up
$ws.Range("H12").Value = 93.375
$ws.Range("J12").Value = 125.818184
$ws.Range("L12").Value = 377.454552
$ws.Range("N12").Value = -723.454552

$ws = $wb.Worksheets.Item("GSM")
# Row 46: Burning the Midnight Oil
$ws.Range("H46").Value = 4421
$ws.Range("J46").Value = 4421
$ws.Range("L46").Value = 4421
$ws.Range("N46").Value = -4733

# Row 80: Needs More Prayerbell
$ws.Range("H80").Value = 17611266
$ws.Range("I80").Value = 43920668
$ws.Range("J80").Value = 71665.164
$ws.Range("K80").Value = 43920668
$ws.Range("L80").Value = 71665.164
$ws.Range("M80").Value = -43919670
$ws.Range("N80").Value = -73661.164

# Row 83: With a Noise That Reaches Heaven (L)
$ws.Range("H83").Value = 17611266
$ws.Range("I83").Value = 43920668
$ws.Range("J83").Value = 71665.164
$ws.Range("K83").Value = 219603340
$ws.Range("L83").Value = 358325.82
$ws.Range("M83").Value = -219598348
$ws.Range("N83").Value = -368309.82

# Row 97: If I'd a Koppranickel for Every Time...
$ws.Range("H97").Value = 2505
$ws.Range("I97").Value = 2505
$ws.Range("K97").Value = 2505
$ws.Range("M97").Value = -2009

$ws = $wb.Worksheets.Item("LTW")
# Row 100: Tiger in the Sack
$ws.Range("H100").Value = 3220.375
$ws.Range("I100").Value = 2927.6667
$ws.Range("J100").Value = 3396
$ws.Range("K100").Value = 2927.6667
$ws.Range("L100").Value = 3396
$ws.Range("M100").Value = -2386.6667
$ws.Range("N100").Value = -4478

$ws = $wb.Worksheets.Item("WVR")
# Row 92: Modest Beginnings
$ws.Range("H92").Value = 100550
$ws.Range("J92").Value = 100550
$ws.Range("L92").Value = 100550
$ws.Range("N92").Value = -105542

# Row 106: Cap It Off
$ws.Range("H106").Value = 20000
$ws.Range("J106").Value = 20000
$ws.Range("L106").Value = 20000
$ws.Range("N106").Value = -22524
